$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.081.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.758.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.16%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.50%  '

$ws.Range("E7").Value = '  +0.22%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.111'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -13.56%  '

$ws.Range("E11").Value = '  -0.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.157'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.247.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("E14").Value = '  +2.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.998.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.42%  '

$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.761.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("E18").Value = '  +1.35%  '

$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '360.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.23%  '

$ws.Range("E22").Value = '  +4.09%  '

$ws.Range("E23").Value = '  +0.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("E25").Value = '  +2.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("E27").Value = '  +0.33%  '

$ws.Range("E28").Value = '  +5.33%  '

$ws.Range("E29").Value = '  -1.90%  '

$ws.Range("E30").Value = '  -0.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.44'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.95%  '

$ws.Range("E35").Value = '  +0.17%  '

$ws.Range("E36").Value = '  +1.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '333.53'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.19%  '

$ws.Range("E42").Value = '  +0.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0597'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.62%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("E46").Value = '  +0.66%  '

$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.46%  '

$ws.Range("E49").Value = '  +0.83%  '

$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.64%  '
